$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.717.84'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.645.67'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.15'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.530'
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.30'
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0613'
$ws.Range("E10").Value = '  +0.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0892'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.880.21'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.650.66'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.63'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.718.14'
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.71'
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0724'
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.62'
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.28'
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.09'
$ws.Range("E23").Value = '  +8.42%  '
$ws.Range("E24").Value = '  -3.86%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.12'
$ws.Range("E25").Value = '  +1.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.92'
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("E28").Value = '  +0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.63'
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("E30").Value = '  +0.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0486'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.443.71'
$ws.Range("E33").Value = '  +3.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("E35").Value = '  +2.19%  '
$ws.Range("E36").Value = '  -1.25%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.566'
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.878'
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0166'
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.886'
$ws.Range("E40").Value = '  +12.85%  '
$ws.Range("E41").Value = '  +1.23%  '
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.58'
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("B44").Value = 'mCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.48'
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '67.02'
$ws.Range("E45").Value = '  +4.47%  '
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.789.38'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("E48").Value = '  +5.04%  '
$ws.Range("E49").Value = '  +2.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '85.46'
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0988'
$ws.Range("E51").Value = '  +0.52%  '